$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make sure date-like text stays as plain text (not auto-converted to a date serial),
# then drop back to the default "Normal" cell style so no stray number format sticks around.
$ws.Range("B2:B4").NumberFormat = "@"

# Update row 2 (product "Kompos" -> "Pupuk") with new date + price
$ws.Range("B2").Value = "2023-12-30"
$ws.Range("C2").Value = "Pupuk"
$ws.Range("D2").Value = 10000

# Update row 3 ("Maggot") with new date, keep price
$ws.Range("B3").Value = "2024-03-18"
$ws.Range("C3").Value = "Maggot"
$ws.Range("D3").Value = 5000

# Add new row 4 - another Maggot setoran (cancel setoran feature add)
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "2024-03-18"
$ws.Range("C4").Value = "Maggot"
$ws.Range("D4").Value = 2000

$ws.Range("B2:B4").Style = "Normal"
